{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"System Monitor\" list-item paragraph and the \"Procedure to pin\"\n// paragraph that immediately follows it (content delivered in the diff).\nlet systemMonitorPara = null;\nlet procedureToPinPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"System Monitor\") {\n    systemMonitorPara = paragraphs.items[i];\n    // the next paragraph should be \"Procedure to pin\"\n    if (i + 1 < paragraphs.items.length) {\n      procedureToPinPara = paragraphs.items[i + 1];\n    }\n    break;\n  }\n}\n\nif (systemMonitorPara) {\n  // The \"_GoBack\" bookmark tracks the last edited location; after removing\n  // the \"System Monitor\" paragraph, Word leaves it at the start of the\n  // paragraph that now follows (\"Procedure to pin\") instead of its old spot\n  // at the end of \"Update all packages and the OS\". Move it there.\n  context.document.deleteBookmark(\"_GoBack\");\n\n  if (procedureToPinPara) {\n    const startRange = procedureToPinPara.getRange(\"Start\");\n    startRange.insertBookmark(\"_GoBack\");\n  }\n\n  systemMonitorPara.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"System Monitor\" list-item paragraph (to be removed) and the\n# paragraph that immediately follows it (\"Procedure to pin\").\n$systemMonitorIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs($i).Range.Text.Trim() -eq \"System Monitor\") {\n    $systemMonitorIndex = $i\n    break\n  }\n}\n\nif ($systemMonitorIndex -ne $null) {\n  # The hidden \"_GoBack\" bookmark marks the location of the last edit. It\n  # currently sits at the end of the \"Update all packages and the OS\"\n  # paragraph; after deleting \"System Monitor\" it needs to move to the start\n  # of the paragraph that now directly follows (\"Procedure to pin\").\n  try {\n    $goBack = $d.Bookmarks.Item(\"_GoBack\")\n    if ($goBack -ne $null) {\n      $goBack.Delete()\n    }\n  } catch {\n    # no pre-existing \"_GoBack\" bookmark - nothing to remove\n  }\n\n  $nextIndex = $systemMonitorIndex + 1\n  if ($nextIndex -le $d.Paragraphs.Count) {\n    $target = $d.Paragraphs($nextIndex).Range.Duplicate()\n    $target.Collapse(1) # wdCollapseStart\n    $target.Bookmarks.Add(\"_GoBack\")\n  }\n\n  $d.Paragraphs($systemMonitorIndex).Range.Delete()\n}\n"}
